# Append new ticker rows to the end of the data (Daten aktualisiert am 2023-12-12)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A204").Value = "IMX-USD"
$ws.Range("A205").Value = "TAO-USD"
$ws.Range("A206").Value = "GRT-USD"
